$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Convert a handful of cells that were numeric typos to the text "NaN" ---
# (these reuse the existing shared string "NaN" already present in the workbook)
$nanCells = @("DT9", "CL18", "CC43", "O51", "DK54", "AP83", "CF109")
foreach ($addr in $nanCells) {
    $ws.Range($addr).Value = "NaN"
}

# --- Corrected figures in column AD for a handful of existing rows ---
$ws.Range("AD163").Value = 65
$ws.Range("AD164").Value = 72
$ws.Range("AD165").Value = 82
$ws.Range("AD166").Value = 99

# --- Append a new day of data as row 190 ---
$row190 = @(44084,694664,2721,92177,65525,234514,26658,4718,3671,7179,6922,14405,3841,21641,27183,6328,6948,13500,10805,15751,13247,3314,1779,7634,23286,12824,9030,52085,1404,430,512,457,337,222,426,1995,4061,36930,7953,2484,40739,1033,21446,1481,9157,1590,1583,5926,1743,953,2484,2652,53341,13172,4008,8497,5089,280,1429,2645,735,2104,9000,8987,9513,14046,1916,855,10800,8638,10133,1927,1770,4192,4023,1278,5140,2865,1611,820,2498,2105,1582,1214,5827,1758,1288,1488,1874,1782,2125,1335,1161,1155,710,3123,1224,842,876,1544,1399,706,804,1090,1377,1190,1303,1017,324,346,747,678,444,534,356,639,724,518,482,372,517,126416,294919,13357,126632,78599,36582,10830)

for ($i = 0; $i -lt $row190.Length; $i++) {
    $ws.Cells.Item(190, $i + 1).Value = $row190[$i]
}

# --- Update the view's active selection to match the latest data row ---
$ws.Range("DS185").Select() | Out-Null
